$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (date, serial number) values between rows 4-5 and 6-7
$ws.Range("D4").Value2 = 44559
$ws.Range("D5").Value2 = 44559
$ws.Range("D6").Value2 = 44574
$ws.Range("D7").Value2 = 44574
